$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    'stecken',
    'brauchen',
    'schwimmen',
    'saufen',
    'husten',
    'malen',
    'folgen',
    'fallen',
    'scheitern',
    'wirken',
    'fehlen',
    'stoppen',
    'flehen',
    'spüren',
    'bauen',
    'schulden',
    'feiern',
    'segnen',
    'sorgen',
    'münzen',
    'geben',
    'wecken',
    'heilen',
    'graben',
    'schützen',
    'irren',
    'morden',
    'bellen',
    'wundern',
    'schwören',
    'backen',
    'klingeln',
    'leugnen',
    'weichen',
    'schenken',
    'tropfen',
    'hören',
    'äußern',
    'dringen',
    'schalten',
    'posten',
    'spinnen',
    'lesen',
    'zeigen',
    'bitten',
    'sichern',
    'räumen',
    'spenden',
    'zögern',
    'formen',
    'bremsen',
    'greifen',
    'töten',
    'warnen',
    'sprengen',
    'dienen',
    'lassen',
    'rauben',
    'boxen',
    'fischen',
    'treiben',
    'schaden',
    'danken',
    'wehtun',
    'mögen',
    'klingen',
    'spielen',
    'beten',
    'sinken',
    'heulen',
    'knarren',
    'schätzen',
    'schrecken',
    'jagen',
    'lügen',
    'schreiten',
    'starten',
    'kennen',
    'stehlen',
    'mühen',
    'führen',
    'jubeln',
    'nennen',
    'enden',
    'weinen',
    'arten',
    'zünden',
    'ändern',
    'sparen',
    'wenden',
    'hacken',
    'kichern',
    'mauern',
    'leihen',
    'kosten',
    'liefern',
    'stürzen',
    'klappen',
    'garen',
    'trauen',
    'werfen',
    'pfeifen',
    'fällen',
    'ehren',
    'streichen',
    'schleppen',
    'schwächen',
    'sterben',
    'baden',
    'ahnen',
    'erben',
    'seufzen',
    'wählen',
    'helfen',
    'wüten',
    'drucken',
    'trösten',
    'fassen',
    'biegen',
    'betteln',
    'schlagen',
    'platzen',
    'stammen',
    'scheinen',
    'nerven',
    'rühren',
    'schließen',
    'flüchten',
    'fangen',
    'zielen',
    'kehren',
    'sprechen',
    'rücken',
    'handeln',
    'grüßen',
    'parken',
    'fügen',
    'bergen',
    'fließen',
    'reizen',
    'gründen',
    'fragen',
    'schneiden',
    'machen',
    'proben',
    'decken',
    'runden',
    'merken',
    'gelten',
    'hauen',
    'stellen',
    'streifen',
    'liegen',
    'wachsen',
    'planen',
    'quälen',
    'kümmern',
    'bluten',
    'fordern',
    'loben',
    'bleiben',
    'lockern',
    'trennen',
    'drehen',
    'schlucken',
    'tanzen',
    'doppeln',
    'schmecken',
    'sperren',
    'filmen',
    'pflanzen',
    'schwingen',
    'kratzen',
    'testen',
    'wahren',
    'rufen',
    'tollen',
    'fahren',
    'kranken',
    'klettern',
    'altern',
    'suchen',
    'freuen',
    'werden',
    'narren',
    'siegen',
    'rasen',
    'hassen',
    'achten',
    'ärgern',
    'gleichen',
    'heben'
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $words[$i]
}

Write-Output "Updated $($words.Length) cells in column B"
